$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# Update existing row 237 (B and D values changed)
$ws.Cells.Item(237, 2).Value = 177772000000000
$ws.Cells.Item(237, 4).Value = 136015302218.8217

# Copy the style of A237 (date-formatted column) for reuse on new rows
$ws.Cells.Item(237, 1).Copy() | Out-Null

# Add new row 238
$ws.Cells.Item(238, 1).PasteSpecial($xlPasteFormats) | Out-Null
$ws.Cells.Item(238, 1).Value = 45108
$ws.Cells.Item(238, 2).Value = 176788000000000
$ws.Cells.Item(238, 3).Value = 0.0007645259938837921
$ws.Cells.Item(238, 4).Value = 135159021406.7278

# Add new row 239
$ws.Cells.Item(239, 1).PasteSpecial($xlPasteFormats) | Out-Null
$ws.Cells.Item(239, 1).Value = 45139
$ws.Cells.Item(239, 2).Value = 174322000000000
$ws.Cells.Item(239, 3).Value = 0.0007641637755803825
$ws.Cells.Item(239, 4).Value = 133210557686.7234

# Add new row 240
$ws.Cells.Item(240, 1).PasteSpecial($xlPasteFormats) | Out-Null
$ws.Cells.Item(240, 1).Value = 45170
$ws.Cells.Item(240, 2).Value = 173950000000000
$ws.Cells.Item(240, 3).Value = 0.0007641637755803825
$ws.Cells.Item(240, 4).Value = 132926288762.2075

$excel.CutCopyMode = 0
